$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date Prepared (C8): 10/9/2024 -> 1/6/2025
$ws.Range("C8").Value = "1/6/2025"

# Completion Date (I8): 10/10/2024 -> 1/6/2024
$ws.Range("I8").Value = "1/6/2024"

# Delivery Date (I9): 10/11/2024 -> 1/7/2024
$ws.Range("I9").Value = "1/7/2024"

# Urgency No. (I10): empty -> 1, picking up the date-style formatting
# used by the neighboring Completion/Delivery Date cells (matches source edit).
$ws.Range("I9").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 1

# JO No. (C10): "ADM-137-2024" -> "ADM-137-20241"
$ws.Range("C10").Value = "ADM-137-20241"

# Restore selection to C10:E10 (also scrolls the view back to the top,
# clearing the previous topLeftCell="A12" scroll position).
[void]$ws.Range("C10:E10").Select()
